$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D10").Value = 500
$ws.Range("D11").Value = 500
$ws.Range("E10").Value = "Internal Server Error"
$ws.Range("E11").Value = "Internal Server Error"
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "teacher_abai_school"
$ws.Range("C12").Value = "=pkFg4tf60AW"
$ws.Range("D12").Value = 200
